# Auto-generated Word COM-interop script implementing the target diff.
# Strategy: for each target span, delete the old (possibly multi-run) text via
# Find.Execute(..., "", wdReplaceAll) which collapses the range to the deletion
# point, then InsertAfter() each desired new run's text in sequence at computed
# offsets -- this creates fresh, separate <w:r> elements (confirmed empirically)
# rather than merging into a neighbouring run.

$d = $word.ActiveDocument

function Replace-Span($OldText, $NewParts) {
    $r = $d.Content
    $found = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
    if (-not $found) {
        throw ("Span not found: " + $OldText)
    }
    $pos = $r.Start
    foreach ($part in $NewParts) {
        $ins = $d.Range($pos, $pos)
        $ins.InsertAfter($part)
        $pos = $pos + $part.Length
    }
}

# Paragraph 1 (Reusability) -- merge 5 runs into 1, text unchanged
$old = 'From these hardware selections it allows our project to be utilized on many different platforms maximizing its reusability. Finally, to further the reusability of our project the TurtleBot’s code is run through a python script that can be converted into an executable file allowing any OS to drive the autonomous robotic code through ROS.'
$newParts = @('From these hardware selections it allows our project to be utilized on many different platforms maximizing its reusability. Finally, to further the reusability of our project the TurtleBot’s code is run through a python script that can be converted into an executable file allowing any OS to drive the autonomous robotic code through ROS.')
Replace-Span $old $newParts

# Paragraph 2 (Reconfigurations) part A -- merge runs around '360 degree' proofErr, text unchanged
$old = 'As a direct consequence of the emphasis placed on reusability in this product, reconfiguration options are both simplistic and plentiful. To begin, the system’s inputs can be simply modified, as the ultrasonic sensors can be replaced, or added to, by other sensors with UART or MQTT communication features. Adding to the input sensors would allow more complicated ranges and complexity of gestures to be realised, with mappings of 360 degree gestures possible through using LiDAR sensors, more intricate gestures can be mapped utilising an input camera and machine learning techniques and so on.'
$newParts = @('As a direct consequence of the emphasis placed on reusability in this product, reconfiguration options are both simplistic and plentiful. To begin, the system’s inputs can be simply modified, as the ultrasonic sensors can be replaced, or added to, by other sensors with UART or MQTT communication features. Adding to the input sensors would allow more complicated ranges and complexity of gestures to be realised, with mappings of 360 degree gestures possible through using LiDAR sensors, more intricate gestures can be mapped utilising an input camera and machine learning techniques and so on.')
Replace-Span $old $newParts

# Paragraph 2 (Reconfigurations) part B -- merge 'Another reconfiguration...' runs, text unchanged
$old = 'Another reconfiguration is by finding a programmable board that supports MQTT communication allowing for the external m5core2 communication node to be disconnected. This cross-communication removal allows for the tracked position of the TurtleBot’s location to be displayed on many different digital display devices that utilise the UART communication protocol instead of the m5core2. This would significantly simplify the overall system, as it removes one of the programmable modules involved and replaces it with a simple UART LCD display. For further features the programmable board used could also have other features such as '
$newParts = @('Another reconfiguration is by finding a programmable board that supports MQTT communication allowing for the external m5core2 communication node to be disconnected. This cross-communication removal allows for the tracked position of the TurtleBot’s location to be displayed on many different digital display devices that utilise the UART communication protocol instead of the m5core2. This would significantly simplify the overall system, as it removes one of the programmable modules involved and replaces it with a simple UART LCD display. For further features the programmable board used could also have other features such as ')
Replace-Span $old $newParts

# Paragraph 2 (Reconfigurations) part C -- merge Bluetooth runs, text unchanged
$old = 'being able to utilise low energy Bluetooth capabilities, which could promote communication to multiple different autonomous robot systems, or to communicate with ranging sensors to further guide the robot. '
$newParts = @('being able to utilise low energy Bluetooth capabilities, which could promote communication to multiple different autonomous robot systems, or to communicate with ranging sensors to further guide the robot. ')
Replace-Span $old $newParts

# Paragraph 3 (Ease of use) -- replace garbled '***\n* .' with new explanatory sentence, split across 5 runs
$old = 'To begin, the TurtleBot’s python script can be converted into a python executable file, **** .'
$newParts = @(
    'To begin, the TurtleBot’s python script can be converted into a python executable file, ',
    'which will ',
    'automatically connect to ROS and run the operational script',
    '.',
    ' This process will enhance efficiency and simplicity, as the deployment team will not be required to access the ROS terminal in order to set up the product.'
)
Replace-Span $old $newParts

# Paragraph 4 (manual/wiring) -- merge all runs into 1, text unchanged
$old = 'Additionally, a manual has been created to describe the gesture actions performed to the Ultrasonic sensors, and what the TurtleBot’s responses to these gestures will be. This allows a deployment team to prepare their actions ahead of deployment, in order to seamlessly deploy the product. Further, despite the intentional simplicity of wiring this device, a wiring manual has been created, which will further aid the deployment team in efficiently setting up the product.'
$newParts = @('Additionally, a manual has been created to describe the gesture actions performed to the Ultrasonic sensors, and what the TurtleBot’s responses to these gestures will be. This allows a deployment team to prepare their actions ahead of deployment, in order to seamlessly deploy the product. Further, despite the intentional simplicity of wiring this device, a wiring manual has been created, which will further aid the deployment team in efficiently setting up the product.')
Replace-Span $old $newParts

Write-Output "All replacements applied."
